# Refresh the cryptos list (Price / Volume(1h) columns) with the latest
# scrape, including the ImmutableX / NEARProtocol row-order swap.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Several 'Price' values look like plain numbers (e.g. 608.72, 0.483).
# Assigning such a string straight to .Value lets Excel auto-convert it
# to a Double, which both changes its type and can introduce binary
# floating-point noise (e.g. 608.72 -> 608.72000000000003). Prefixing
# with an apostrophe is the normal Excel way to force literal text,
# matching the original inline-string cell content exactly.
function Set-TextValue($range, $value) {
    if ($value -match '^\s*[+-]?\d+(\.\d+)?\s*$') {
        $range.Value = "'" + $value
    } else {
        $range.Value = $value
    }
}

Set-TextValue $ws.Range('D2') '67.099.24'
Set-TextValue $ws.Range('E2') '  -3.30%  '
Set-TextValue $ws.Range('D3') '3.531.82'
Set-TextValue $ws.Range('E3') '  -3.87%  '
Set-TextValue $ws.Range('E4') '  +0.16%  '
Set-TextValue $ws.Range('D5') '608.72'
Set-TextValue $ws.Range('E5') '  -5.69%  '
Set-TextValue $ws.Range('D6') '152.90'
Set-TextValue $ws.Range('E6') '  -3.53%  '
Set-TextValue $ws.Range('D7') '3.528.15'
Set-TextValue $ws.Range('E7') '  -3.98%  '
Set-TextValue $ws.Range('E8') '  +0.16%  '
Set-TextValue $ws.Range('D9') '0.483'
Set-TextValue $ws.Range('E9') '  -2.95%  '
Set-TextValue $ws.Range('E10') '  -3.22%  '
Set-TextValue $ws.Range('D11') '6.82'
Set-TextValue $ws.Range('E11') '  -3.85%  '
Set-TextValue $ws.Range('D12') '0.427'
Set-TextValue $ws.Range('E12') '  -3.81%  '
Set-TextValue $ws.Range('E13') '  -3.83%  '
Set-TextValue $ws.Range('D14') '4.130.38'
Set-TextValue $ws.Range('E14') '  -3.74%  '
Set-TextValue $ws.Range('D15') '31.74'
Set-TextValue $ws.Range('E15') '  -2.27%  '
Set-TextValue $ws.Range('D16') '3.527.43'
Set-TextValue $ws.Range('E16') '  -4.36%  '
Set-TextValue $ws.Range('D17') '67.109.89'
Set-TextValue $ws.Range('E17') '  -3.25%  '
Set-TextValue $ws.Range('E18') '  +1.06%  '
Set-TextValue $ws.Range('D19') '6.32'
Set-TextValue $ws.Range('E19') '  -2.25%  '
Set-TextValue $ws.Range('D20') '15.39'
Set-TextValue $ws.Range('E20') '  -2.74%  '
Set-TextValue $ws.Range('D21') '444.23'
Set-TextValue $ws.Range('E21') '  -5.01%  '
Set-TextValue $ws.Range('D22') '9.24'
Set-TextValue $ws.Range('E22') '  -7.65%  '
Set-TextValue $ws.Range('D23') '0.628'
Set-TextValue $ws.Range('E23') '  -2.50%  '
Set-TextValue $ws.Range('D24') '77.79'
Set-TextValue $ws.Range('E24') '  -1.89%  '
Set-TextValue $ws.Range('D25') '3.676.28'
Set-TextValue $ws.Range('E25') '  -3.71%  '
Set-TextValue $ws.Range('E26') '  +0.00%  '
Set-TextValue $ws.Range('E27') '  -2.08%  '
Set-TextValue $ws.Range('D28') '10.20'
Set-TextValue $ws.Range('E28') '  -4.61%  '
Set-TextValue $ws.Range('D29') '8.25'
Set-TextValue $ws.Range('E29') '  -7.95%  '
Set-TextValue $ws.Range('D30') '2.53'
Set-TextValue $ws.Range('E30') '  -3.37%  '
Set-TextValue $ws.Range('D31') '1.67'
Set-TextValue $ws.Range('E31') '  -0.85%  '
Set-TextValue $ws.Range('E32') '  +0.14%  '
Set-TextValue $ws.Range('D33') '25.68'
Set-TextValue $ws.Range('E33') '  -4.73%  '
Set-TextValue $ws.Range('D34') '0.158'
Set-TextValue $ws.Range('E34') '  -2.13%  '
Set-TextValue $ws.Range('B35') 'ImmutableX'
Set-TextValue $ws.Range('C35') 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range('D35') '1.87'
Set-TextValue $ws.Range('E35') '  -6.55%  '
Set-TextValue $ws.Range('B36') 'NEARProtocol'
Set-TextValue $ws.Range('C36') 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range('D36') '6.14'
Set-TextValue $ws.Range('E36') '  -4.06%  '
Set-TextValue $ws.Range('D37') '3.527.27'
Set-TextValue $ws.Range('E37') '  -3.75%  '
Set-TextValue $ws.Range('D38') '7.99'
Set-TextValue $ws.Range('E38') '  -4.67%  '
Set-TextValue $ws.Range('E40') '  +0.11%  '
Set-TextValue $ws.Range('D41') '175.69'
Set-TextValue $ws.Range('E41') '  -1.68%  '
Set-TextValue $ws.Range('E42') '  -2.74%  '
Set-TextValue $ws.Range('D43') '5.55'
Set-TextValue $ws.Range('E43') '  -5.23%  '
Set-TextValue $ws.Range('D44') '0.0860'
Set-TextValue $ws.Range('E44') '  -3.39%  '
Set-TextValue $ws.Range('D45') '0.890'
Set-TextValue $ws.Range('E45') '  -3.81%  '
Set-TextValue $ws.Range('D46') '45.62'
Set-TextValue $ws.Range('E46') '  -2.76%  '
Set-TextValue $ws.Range('D47') '27.47'
Set-TextValue $ws.Range('E47') '  -2.26%  '
Set-TextValue $ws.Range('D48') '2.60'
Set-TextValue $ws.Range('E48') '  -3.20%  '
Set-TextValue $ws.Range('E49') '  -1.44%  '
Set-TextValue $ws.Range('E50') '  -2.93%  '
Set-TextValue $ws.Range('D51') '7.56'
Set-TextValue $ws.Range('E51') '  -2.87%  '
